$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TASK SUMMARY SHEET")

$ws.Range("C1").Value = "Jesse Hare"
$ws.Range("E1").Value = 5

$ws.Range("A14").Value = "Cumulative Total: 120"

$ws.Range("A3").Value = "Project Build"
$ws.Range("B3").Value = "Work on first iteration to present to client"
$ws.Range("C3").Value = 20
$ws.Range("D3").Value = 20

$ws.Range("D14").Value = 20

[void]$ws.Range("D11").Select()
